$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells (F1, G1) and updated / new data cells.
#    Written in this order so the shared-string table is populated in the
#    same sequence as the target workbook (Country code, new D2 email,
#    new D3 email, Address, Devas Naka/Indore, Pune).
# ---------------------------------------------------------------------------
$ws.Range("F1").Value = "Country code"
$ws.Range("D2").Value = "mayur.charvande+4@fxbytes1.com"
$ws.Range("D3").Value = "mayur.charvande+3@fxbytes2.com"
$ws.Range("G1").Value = "Address"
$ws.Range("G2").Value = "Devas Naka, Indore"
$ws.Range("G3").Value = "Pune"
$ws.Range("F2").Value = 91
$ws.Range("F3").Value = 92

# ---------------------------------------------------------------------------
# 2. Column widths for the two new columns.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 18
$ws.Columns("G").ColumnWidth = 27

# ---------------------------------------------------------------------------
# 3. Left-align every populated cell (header row, plain data cells, and the
#    hyperlinked e-mail cells each keep their own font/fill, just gaining a
#    left horizontal alignment).
# ---------------------------------------------------------------------------
$ws.Range("A1:G1").HorizontalAlignment = -4131
$ws.Range("A2:B3").HorizontalAlignment = -4131
$ws.Range("E2:G3").HorizontalAlignment = -4131
$ws.Range("C2:D3").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# 4. Selection / active cell, matching the saved view state.
# ---------------------------------------------------------------------------
$ws.Range("F9").Select()

# ---------------------------------------------------------------------------
# 5. Page setup (printer paper size + orientation).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "edit applied"
